$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new log entry row (row 9), reusing the date style from the row above
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9").Value = 43544
$ws.Range("B9").Value = "Completed all three apis, fixed errors and recorded the demo video"

# Update selection to match the new active cell (B9), matching the diff
$ws.Range("B9").Select()
